# Generate Report for Handback
# Adds two new handed-back files to the Overview / zh-cn / de-de sheets:
#   1f48125f-c151-49d2-9332-e774cba1ddb4.md
#   a4835ca2-106e-4658-a384-828548b54b5f.md

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$statusText = "Handed back: in sync with en-US"
$includeText = "Include"

# ---- File 1 : 1f48125f-c151-49d2-9332-e774cba1ddb4 -----------------------
$f1Name       = "1f48125f-c151-49d2-9332-e774cba1ddb4.md"
$f1XlfZh      = "1f48125f-c151-49d2-9332-e774cba1ddb4.4c2932a2761aea97599b1f8af3b2494fff007153.zh-cn.xlf"
$f1XlfDe      = "1f48125f-c151-49d2-9332-e774cba1ddb4.4c2932a2761aea97599b1f8af3b2494fff007153.de-de.xlf"
$f1HandoffDtZh  = "2016-02-22 05:17:45"
$f1HandbackDtZh = "2016-02-22 05:18:48"
$f1HandoffDtDe  = "2016-02-22 05:17:58"
$f1HandbackDtDe = "2016-02-22 05:19:11"

$f1UrlSrc         = "https://github.com/OpenLocalizationTest/oltest/blob/06a99f80d9d76bae16507416e340e09ecc710a6d/e2e/$f1Name"
$f1UrlHandoffZh   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0076ce2e351023f1e45978a685901f18f801c1d6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$f1XlfZh"
$f1UrlTargetZh    = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0e70b407e14ee7965af5504b3b4b008d234f4cbc/e2e/$f1Name"
$f1UrlHandbackZh  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ace6424d75d594dac4cc917f70ee8219cc15e510/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$f1XlfZh"
$f1UrlHandoffDe   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4360603698904b931e2d0f872734ebb28cd540c7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$f1XlfDe"
$f1UrlTargetDe    = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/baff15f3c13079d8d3e2e1a8d713fe6c77d8fdd0/e2e/$f1Name"
$f1UrlHandbackDe  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9c81aee1c1479c88b76c9bd98ac9456c52b9575e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$f1XlfDe"

# ---- File 2 : a4835ca2-106e-4658-a384-828548b54b5f -----------------------
$f2Name       = "a4835ca2-106e-4658-a384-828548b54b5f.md"
$f2XlfZh      = "a4835ca2-106e-4658-a384-828548b54b5f.68670cf37816f4cd7715229f0f9c091122526fd5.zh-cn.xlf"
$f2XlfDe      = "a4835ca2-106e-4658-a384-828548b54b5f.68670cf37816f4cd7715229f0f9c091122526fd5.de-de.xlf"
$f2HandoffDtZh  = "2016-02-22 05:17:45"
$f2HandbackDtZh = "2016-02-22 05:18:48"
$f2HandoffDtDe  = "2016-02-22 05:17:58"
$f2HandbackDtDe = "2016-02-22 05:19:11"

$f2UrlSrc         = "https://github.com/OpenLocalizationTest/oltest/blob/afd4bdd1ffc09230ebd31ee7ffed9aa62c95305b/e2e/$f2Name"
$f2UrlHandoffZh   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9361d117429eca464c64ba7ddd38fc82322b590/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$f2XlfZh"
$f2UrlTargetZh    = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7b57fcc3efdb04313cd9a9537cf8b5b8200eecb9/e2e/$f2Name"
$f2UrlHandbackZh  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cbb3229d757e38bcafbeddf700d62df19ac3f353/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$f2XlfZh"
$f2UrlHandoffDe   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da53671402fb1b643415bb17383f728f21a1816a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$f2XlfDe"
$f2UrlTargetDe    = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/00047918020caaad12eb0e72cf3346a676b880cd/e2e/$f2Name"
$f2UrlHandbackDe  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e3ab56d91af89606aa0ecbadd089187e83b17c24/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$f2XlfDe"

function Set-LinkCell($range, $text, $url) {
    $range.Value = $text
    $range.Font.Underline = $true
    $range.Font.Color = 0x6495ED
    $range.Worksheet.Hyperlinks.Add($range, $url, "", "", $text) | Out-Null
}

function Set-DateCell($range, $text) {
    $range.Value = $text
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# =====================================================================
# Sheet 1 : Overview -- rows 6 and 7
# =====================================================================
Set-LinkCell $ws1.Range("A6") $f1Name $f1UrlSrc
$ws1.Range("B6").Value = $statusText
$ws1.Range("C6").Value = $statusText

Set-LinkCell $ws1.Range("A7") $f2Name $f2UrlSrc
$ws1.Range("B7").Value = $statusText
$ws1.Range("C7").Value = $statusText

# =====================================================================
# Sheet 2 : zh-cn -- rows 6 and 7
# =====================================================================
Set-LinkCell $ws2.Range("A6") $f1Name $f1UrlSrc
$ws2.Range("B6").Value = $statusText
Set-LinkCell $ws2.Range("C6") $f1XlfZh $f1UrlHandoffZh
Set-DateCell $ws2.Range("D6") $f1HandoffDtZh
Set-LinkCell $ws2.Range("E6") $f1Name $f1UrlTargetZh
Set-LinkCell $ws2.Range("F6") $f1XlfZh $f1UrlHandbackZh
Set-DateCell $ws2.Range("G6") $f1HandbackDtZh
$ws2.Range("H6").Value = $includeText

Set-LinkCell $ws2.Range("A7") $f2Name $f2UrlSrc
$ws2.Range("B7").Value = $statusText
Set-LinkCell $ws2.Range("C7") $f2XlfZh $f2UrlHandoffZh
Set-DateCell $ws2.Range("D7") $f2HandoffDtZh
Set-LinkCell $ws2.Range("E7") $f2Name $f2UrlTargetZh
Set-LinkCell $ws2.Range("F7") $f2XlfZh $f2UrlHandbackZh
Set-DateCell $ws2.Range("G7") $f2HandbackDtZh
$ws2.Range("H7").Value = $includeText

# =====================================================================
# Sheet 3 : de-de -- rows 6 and 7
# =====================================================================
Set-LinkCell $ws3.Range("A6") $f1Name $f1UrlSrc
$ws3.Range("B6").Value = $statusText
Set-LinkCell $ws3.Range("C6") $f1XlfDe $f1UrlHandoffDe
Set-DateCell $ws3.Range("D6") $f1HandoffDtDe
Set-LinkCell $ws3.Range("E6") $f1Name $f1UrlTargetDe
Set-LinkCell $ws3.Range("F6") $f1XlfDe $f1UrlHandbackDe
Set-DateCell $ws3.Range("G6") $f1HandbackDtDe
$ws3.Range("H6").Value = $includeText

Set-LinkCell $ws3.Range("A7") $f2Name $f2UrlSrc
$ws3.Range("B7").Value = $statusText
Set-LinkCell $ws3.Range("C7") $f2XlfDe $f2UrlHandoffDe
Set-DateCell $ws3.Range("D7") $f2HandoffDtDe
Set-LinkCell $ws3.Range("E7") $f2Name $f2UrlTargetDe
Set-LinkCell $ws3.Range("F7") $f2XlfDe $f2UrlHandbackDe
Set-DateCell $ws3.Range("G7") $f2HandbackDtDe
$ws3.Range("H7").Value = $includeText
